$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.618.36'
$ws.Range("E2").Value = '  +1.96%  '

$ws.Range("D3").Value = '3.816.15'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("D4").Formula = '''1.00'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Formula = '''668.83'
$ws.Range("E5").Value = '  +7.60%  '

$ws.Range("D6").Formula = '''169.23'
$ws.Range("E6").Value = '  +2.42%  '

$ws.Range("D7").Value = '3.814.31'
$ws.Range("E7").Value = '  +1.03%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Formula = '''0.528'
$ws.Range("E9").Value = '  +1.36%  '

$ws.Range("E10").Value = '  +0.37%  '

$ws.Range("D11").Formula = '''0.464'
$ws.Range("E11").Value = '  +2.71%  '

$ws.Range("D12").Formula = '''7.08'
$ws.Range("E12").Value = '  +6.45%  '

$ws.Range("D13").Formula = '''0.0000246'
$ws.Range("E13").Value = '  -0.82%  '

$ws.Range("D14").Formula = '''35.96'
$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("D15").Value = '4.457.10'
$ws.Range("E15").Value = '  +1.06%  '

$ws.Range("D16").Value = '3.808.84'
$ws.Range("E16").Value = '  +1.61%  '

$ws.Range("D17").Value = '70.516.63'
$ws.Range("E17").Value = '  +1.82%  '

$ws.Range("D18").Formula = '''17.75'
$ws.Range("E18").Value = '  +0.56%  '

$ws.Range("D19").Formula = '''7.17'
$ws.Range("E19").Value = '  +1.02%  '

$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").Formula = '''10.88'
$ws.Range("E21").Value = '  +13.08%  '

$ws.Range("D22").Formula = '''475.51'
$ws.Range("E22").Value = '  +1.58%  '

$ws.Range("D23").Formula = '''0.715'
$ws.Range("E23").Value = '  +1.77%  '

$ws.Range("D24").Formula = '''83.02'
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("D25").Formula = '''0.0000145'
$ws.Range("E25").Value = '  -3.53%  '

$ws.Range("D26").Formula = '''12.26'
$ws.Range("E26").Value = '  +1.95%  '

$ws.Range("D27").Formula = '''10.37'
$ws.Range("E27").Value = '  +3.23%  '

$ws.Range("E28").Value = '  -1.55%  '

$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("D30").Value = '3.964.31'
$ws.Range("E30").Value = '  +0.98%  '

$ws.Range("D31").Formula = '''2.87'
$ws.Range("E31").Value = '  +8.15%  '

$ws.Range("D33").Formula = '''7.40'
$ws.Range("E33").Value = '  +1.32%  '

$ws.Range("D34").Formula = '''29.68'
$ws.Range("E34").Value = '  +3.09%  '

$ws.Range("E35").Value = '  +10.26%  '

$ws.Range("D36").Formula = '''9.18'
$ws.Range("E36").Value = '  +2.18%  '

$ws.Range("D37").Formula = '''0.998'
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Formula = '''0.102'
$ws.Range("E38").Value = '  +0.51%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Formula = '''3.47'
$ws.Range("E39").Value = '  +1.85%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Formula = '''5.98'
$ws.Range("E40").Value = '  +2.97%  '

$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").Formula = '''0.969'
$ws.Range("E41").Value = '  +0.19%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Formula = '''1.00'
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Formula = '''2.09'
$ws.Range("E44").Value = '  +9.76%  '

$ws.Range("B45").Value = 'Arweave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D45").Formula = '''45.69'
$ws.Range("E45").Value = '  +5.23%  '

$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Formula = '''156.91'
$ws.Range("E46").Value = '  +2.42%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Formula = '''48.11'
$ws.Range("E47").Value = '  +2.85%  '

$ws.Range("B48").Value = 'TheGraph'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D48").Formula = '''0.301'
$ws.Range("E48").Value = '  +0.62%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Formula = '''1.42'
$ws.Range("E49").Value = '  +3.80%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Formula = '''0.000290'
$ws.Range("E50").Value = '  +4.06%  '

$ws.Range("D51").Formula = '''8.52'
$ws.Range("E51").Value = '  +1.16%  '
